$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (before current row 2),
# pushing the existing 20 data rows down to rows 5-24.
$insertRange = $ws.Range("A2:C4")
$insertRange.EntireRow.Insert()
# The insert copies formatting from the row above (the header); clear it so
# the new data rows match the unstyled look of the rest of the data rows.
$ws.Range("A2:C4").ClearFormats()

# New data for the 3 inserted rows (now rows 2-4)
$newTop = @(
    @(-0.2407464981079101, 0.5433270186185837, -0.3658644706010817),
    @(-0.2030138969421392, 0.5583634674549103, -0.1664191037416463),
    @(-0.3572314977645862, 0.5254133790731429, -0.5817861706018449)
)

for ($i = 0; $i -lt $newTop.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# Append 7 new rows at the end (rows 25-31)
$newBottom = @(
    @(-0.0122048854827883, 0.3785421848297117, -0.2985985279083258),
    @(0.5578445792198222, 0.3809743523597721, 0.0945302546024358),
    @(1.200664520263672, 0.5283758044242859, 0.6270142197608946),
    @(0.02087068557739279, 0.6865898966789239, -0.5058017373085018),
    @(0.09223079681396533, 0.5110847949981689, -0.3808159828186036),
    @(0.2094589471817022, 0.502252608537674, -0.392595499753952),
    @(0.01102042198181016, 0.4913336634635925, -0.3189654350280756)
)

for ($i = 0; $i -lt $newBottom.Count; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value = $newBottom[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottom[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottom[$i][2]
}
